$d = $word.ActiveDocument

# --- Edit 1: Expand the "advanced evaluation function" discussion ---
$old1 = 'I originally had a more advanced evaluation function, where I included an extra evaluation to account for the fact that it is better to have your pieces further advanced along the board, however, after testing, I found that this made my program worse in efficiency and optimality by a fairly large margin, so I decided to omit it from the final submission. I encountered the same efficiency issue when trying to implement an evaluation that assesses the mobility of king pieces, as well as how threatened the king piece is by determining the number of opposing pieces that were diagonal to it.'
$new1 = 'I included an extra evaluation to account for the fact that it is better to have your pieces further advanced along the board, specifically by adding (for red) or subtracting (for black) the row position of each piece to the evaluation. I also implemented an evaluation that assesses the mobility of king pieces, as well as how threatened the king piece is by determining the number of opposing pieces that were diagonal to it. For each opposing piece that threatened the king, 3 points were added (for black) or subtracted (for red) from the evaluation depending on the colour. Moreover, for each available move for the king piece, a value of 1 was added (for red) or subtracted (for black) from the evaluation. This feature coincides with the idea that a player wants to maximize the number of king pieces that they have as they play the game.'
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Edit 2: Remove the now-obsolete "For these reasons..." sentence ---
$old2 = 'For these reasons, these advanced features were omitted from my evaluation function. '
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Edit 3: Collapse the two trailing empty paragraphs into one, with
#     simplified formatting (no paragraph-mark run formatting, no first-line indent) ---
$countBefore = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($countBefore - 1)
$secondLast.Range.Delete()

$countAfter = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($countAfter)
$last.Range.Style = "Normal"
$last.LeftIndent = 0

Write-Output "paragraphs before=$countBefore after=$countAfter"
